# Update the "mapa_interactivo" workbook:
#  1. Column D ("Comuna") values that were stored as long text
#     ("106594 - PALERMO") are re-written as the short numeric Comuna code
#     ("14") across every worksheet.
#  2. The case "6221" (FERNANDEZ 1549) row is removed from the "General"
#     sheet and from the "AYKO" sheet (its filtered copy) - everything below
#     it shifts up by one row.

$wb = $excel.ActiveWorkbook

# Comuna long-name -> numeric code lookup (built from the observed data).
$comunaMap = @{
    "106552 - CONSTITUCION"      = "1"
    "106560 - RECOLETA"          = "2"
    "106556 - SAN CRISTOBAL"     = "3"
    "106558 - BALVANERA"         = "3"
    "106557 - PQUE. PATRICIOS"   = "4"
    "106562 - BOEDO"             = "5"
    "106568 - ALMAGRO"           = "5"
    "106569 - CABALLITO"         = "6"
    "106576 - LINIERS"           = "9"
    "106573 - FLORESTA"          = "10"
    "106580 - SANTA RITA"        = "11"
    "106584 - VILLA GRAL. MITRE" = "11"
    "106589 - VILLA DEVOTO"      = "11"
    "106590 - COGHLAN"           = "12"
    "106591 - VILLA PUEYRREDON"  = "12"
    "106592 - VILLA URQUIZA"     = "12"
    "106593 - SAAVEDRA"          = "12"
    "106581 - COLEGIALES"        = "13"
    "106595 - BELGRANO"          = "13"
    "106596 - NUNEZ"             = "13"
    "106594 - PALERMO"           = "14"
    "106582 - VILLA CRESPO"      = "15"
    "106586 - VILLA ORTUZAR"     = "15"
    "106587 - AGRONOMIA"         = "15"
    "106597 - PARQUE CHAS"       = "15"
}

# 1) Remove the "6221" / FERNANDEZ 1549 record from General + AYKO first
#    (so the comuna rewrite below walks the already-shifted rows).
$general = $wb.Worksheets.Item("General")
$ayko = $wb.Worksheets.Item("AYKO")

for ($r = 1; $r -le $general.UsedRange.Rows.Count; $r++) {
    if ($general.Cells.Item($r, 1).Value2 -eq "6221" -and $general.Cells.Item($r, 3).Value2 -eq "FERNANDEZ 1549") {
        $general.Rows.Item($r).Delete() | Out-Null
        break
    }
}

for ($r = 1; $r -le $ayko.UsedRange.Rows.Count; $r++) {
    if ($ayko.Cells.Item($r, 1).Value2 -eq "6221" -and $ayko.Cells.Item($r, 3).Value2 -eq "FERNANDEZ 1549") {
        $ayko.Rows.Item($r).Delete() | Out-Null
        break
    }
}

# 2) Rewrite column D ("Comuna") on every worksheet: long "code - NAME" text
#    becomes the bare numeric Comuna code.
foreach ($ws in $wb.Worksheets) {
    $lastRow = $ws.UsedRange.Rows.Count
    for ($r = 1; $r -le $lastRow; $r++) {
        $cell = $ws.Cells.Item($r, 4)
        $v = $cell.Value2
        if ($v -ne $null -and $comunaMap.ContainsKey($v)) {
            $cell.Value = $comunaMap[$v]
        }
    }
}

"Done"
